# Code cleanup of the genres_list sheet.
# - "Post-Punk/Darkwave" genre renamed to "Goth / Darkwave"
# - "Synthés" family replaced by "Goths & Corbacs" (Dark Synth - Dungeon synth
#   now classified under the Goth family instead of its own Synthés family)
# - "Pop Rock - Indie Rock" re-family'd from "Blues / Rock" to the new
#   "Post-Punk / Indie Rock" family
# - New genre row added: "Post-Punk" -> family "Post-Punk / Indie Rock"
#
# NOTE: order matters here only insofar as it controls the order new
# strings are appended to the shared-string table; it mirrors the order
# in which the author's edit produced new unique strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename genre in row 28.
$ws.Range("A28").Value = "Goth / Darkwave"

# Add the new "Post-Punk" genre row (row 55) - set A55 before B31 so the
# shared string "Post-Punk" is registered ahead of "Post-Punk / Indie Rock".
$ws.Range("A55").Value = "Post-Punk"

# Match the formatting used by the rest of column A (style used by A2:A54).
$ws.Range("A54").Copy()
$ws.Range("A55").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-family "Pop Rock - Indie Rock" (row 31).
$ws.Range("B31").Value = "Post-Punk / Indie Rock"

# Re-family "Dark Synth - Dungeon synth" (row 42).
$ws.Range("B42").Value = "Goths & Corbacs"

# Family for the newly added "Post-Punk" genre row.
$ws.Range("B55").Value = "Post-Punk / Indie Rock"

# Restore the view state (selection / zoom) as left by the author.
$win = $excel.ActiveWindow
$win.Zoom = 235
$ws.Range("B31").Select()
